# Add a new "Consumption" column into the header row of the Data sheet.
# The new column is inserted at column J (10th column), pushing the
# existing Customer/Tier/Marketplace/... columns one slot to the right
# (J:V -> K:W) and widening the autofilter / used range from A1:V1 to A1:W1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J; this shifts J:V -> K:W and extends
# the sheet's dimension/autofilter automatically.
$ws.Columns("J:J").Insert()

# Populate the new header cell.
$ws.Range("J1").Value = "Consumption"

# Match the style/format of the other header cells (copy A1's formatting,
# i.e. the grey header fill, onto the new header cell) without disturbing
# the shared style table.
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Give the new "Consumption" column a sensible width, as in the template.
$ws.Columns("J:J").ColumnWidth = 16.333333333333332

# Refresh the autofilter so it covers the newly widened header range. The
# sheet already has an autofilter on, so toggle it off first (otherwise
# re-invoking AutoFilter() on an active filter just switches it off).
$ws.AutoFilterMode = $false
$ws.Range("A1:W1").AutoFilter() | Out-Null

# The hidden "_FilterDatabase" defined name backing the autofilter isn't
# auto-expanded by the column insert, so point it at the new range too.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$W`$1"
    }
}

# Update the active selection to the newly added header cell.
$ws.Range("J1").Select()
